$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $result = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- Indicator days block: update the figures in the "Total de dias" run ---
# (the whole block lives in one run, text lines separated by manual line breaks "^l")

$old1 = "Concluído (95 dias)^lAssinatura Contrato (104 dias)"
$new1 = "Concluído (97 dias)^lAssinatura Contrato (106 dias)"
Replace-Text $old1 $new1

$old2 = "Sessão Pública (128 dias)^lHomologado (1 dias)"
$new2 = "Sessão Pública (146 dias)"
Replace-Text $old2 $new2

$old3 = "Homologado (8 dias)^lEm recurso (17 dias)^lAssinatura Contrato (81 dias)^lAssinatura Contrato (67 dias)^lHomologado (0 dias)^lAGU (15 dias)^lTotal de dias 698"
$new3 = "Homologado (8 dias)^lHomologado (2 dias)^lAssinatura Contrato (61 dias)^lAssinatura Contrato (14 dias)^lConcluído (22 dias)^lHomologado (21 dias)^lAssinatura Contrato (49 dias)^lTotal de dias 708"
Replace-Text $old3 $new3

# --- Process numbers: renumber from 0XX to 0(XX-17) ---

Replace-Text "785810/2024-027/00" "785810/2024-010/00"
Replace-Text "785810/2024-028/00" "785810/2024-011/00"
Replace-Text "785810/2024-029/00" "785810/2024-012/00"
Replace-Text "785810/2024-030/00" "785810/2024-013/00"
Replace-Text "785810/2024-031/00" "785810/2024-014/00"
Replace-Text "785810/2024-032/00" "785810/2024-015/00"
Replace-Text "785810/2024-033/00" "785810/2024-016/00"
Replace-Text "785810/2024-034/00" "785810/2024-017/00"
Replace-Text "785810/2024-035/00" "785810/2024-018/00"
Replace-Text "785810/2024-036/00" "785810/2024-019/00"
Replace-Text "785810/2024-037/00" "785810/2024-020/00"
Replace-Text "785810/2024-038/00" "785810/2024-021/00"

Write-Output "Done"
